$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "LSTM Model 1"
$ws.Range("B4").Value = "LSTM Model 2"
$ws.Range("B5").Value = "LSTM Model 3"
$ws.Range("B6").Value = "LSTM Model 1*"
